$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column P (year 2022) that mirrors the formatting of column O
# for each row, then fill in the 2022 values.

# Row 3 (bottom border separator row, empty cell)
$ws.Range("O3").Copy()
$ws.Range("P3").PasteSpecial(-4122)

# Row 4 (year header row)
$ws.Range("O4").Copy()
$ws.Range("P4").PasteSpecial(-4122)
$ws.Range("P4").Value = 2022

# Row 5 (sub-header row, empty cell)
$ws.Range("O5").Copy()
$ws.Range("P5").PasteSpecial(-4122)

# Row 6 data
$ws.Range("O6").Copy()
$ws.Range("P6").PasteSpecial(-4122)
$ws.Range("P6").Value = 1373

# Row 7 data - uses a new right-aligned style with a text dash value
$ws.Range("O7").Copy()
$ws.Range("P7").PasteSpecial(-4122)
$ws.Range("P7").Value = "-"
$ws.Range("P7").HorizontalAlignment = -4152

# Row 8 data
$ws.Range("O8").Copy()
$ws.Range("P8").PasteSpecial(-4122)
$ws.Range("P8").Value = 117

# Row 9 data
$ws.Range("O9").Copy()
$ws.Range("P9").PasteSpecial(-4122)
$ws.Range("P9").Value = 154

# Row 10 data
$ws.Range("O10").Copy()
$ws.Range("P10").PasteSpecial(-4122)
$ws.Range("P10").Value = 885

# Update the selected cell to match the saved view state
$ws.Range("P7").Select()
